$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VisioData")
$ws.Range("A1").Value = "test"
